$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Italicize the leading "set name" letter (and, in one case, the
#    following space) in four short formula paragraphs:
#       "C = {Capitales de Colombia}"
#       "O = {Múltiplos de 8}"
#       "M = {Nombres de mujer}"   (letter AND following space run)
#       "S = {Múltiplos de 7}"
# ---------------------------------------------------------------------

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    if ($t -match "^C = \{Capitales de Colombia\}") {
        $r = $p.Range
        $rng = $d.Range($r.Start, $r.Start + 1)
        $rng.Font.Italic = $true
    }
    elseif ($t -match "^O = \{M.ltiplos de 8\}") {
        $r = $p.Range
        $rng = $d.Range($r.Start, $r.Start + 1)
        $rng.Font.Italic = $true
    }
    elseif ($t -match "^M = \{Nombres de mujer\}") {
        $r = $p.Range
        # Italicize "M" plus the following space run.
        $rng = $d.Range($r.Start, $r.Start + 2)
        $rng.Font.Italic = $true
    }
    elseif ($t -match "^S = \{M.ltiplos de 7\}") {
        $r = $p.Range
        $rng = $d.Range($r.Start, $r.Start + 1)
        $rng.Font.Italic = $true
    }
}

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark: it currently sits just before the
#    "{Múltiplos de 7}" run in the "S = {Múltiplos de 7}" paragraph;
#    it needs to move to the empty (bold/red) paragraph that
#    immediately precedes the "Contenedor 4 ... .mp3)" paragraph.
#    Word keeps only a single "_GoBack" bookmark, so re-adding it
#    under that name relocates it (removing the old one) instead of
#    creating a duplicate. A range spanning just the target
#    paragraph's own mark is not enough to anchor the insertion
#    correctly in this engine, so the range is widened to include the
#    end of the previous paragraph as well (harmless — both are empty
#    paragraph marks with no text).
# ---------------------------------------------------------------------

$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    $t = $p.Range.Text
    if ($t -match "^Contenedor 4 \(25 caracteres") {
        # The bold/red empty paragraph immediately precedes this one.
        $targetIdx = $idx - 1
    }
}

$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($idx -eq $targetIdx) {
        $prevEnd = $p.Range.Start
        $rng = $d.Range($prevEnd - 1, $p.Range.End)
        $d.Bookmarks.Add("_GoBack", $rng)
    }
}
